$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "VIMMP_DEF"

for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 6).Value = "[]"
}
